# Add back the fastq sample rows (35-37) that were found in LTS but had
# no metadata files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 35
Set-TextCell 35 1 "08.14.18"
Set-TextCell 35 2 "H.BROWN"
$ws.Cells.Item(35, 3).Value = 35
Set-TextCell 35 4 "08.14.18"
Set-TextCell 35 5 "H.BROWN"
$ws.Cells.Item(35, 6).Value = 35
$ws.Cells.Item(35, 7).Value = "E7420L"

# Row 36
Set-TextCell 36 1 "10.18.18"
Set-TextCell 36 2 "H.BROWN"
$ws.Cells.Item(36, 3).Value = 36
Set-TextCell 36 4 "10.18.18"
Set-TextCell 36 5 "H.BROWN"
$ws.Cells.Item(36, 6).Value = 36
$ws.Cells.Item(36, 7).Value = "E7420L"

# Row 37
Set-TextCell 37 1 "10.18.18"
Set-TextCell 37 2 "H.BROWN"
$ws.Cells.Item(37, 3).Value = 37
Set-TextCell 37 4 "10.18.18"
Set-TextCell 37 5 "H.BROWN"
$ws.Cells.Item(37, 6).Value = 36
$ws.Cells.Item(37, 7).Value = "E7420L"

# Freeze the header row and move the selection down to the new first empty row,
# matching the author's view state after appending the rows.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A38").Select()
